$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data rows (2-10) were reshuffled (re-sorted), carrying along the
# Fecha (D), Volumen (M), Precio minimo (N), Precio maximo (O),
# Precio promedio ponderado (P) and Precio $/Kg (S) values to new rows.
# Rows 5, 9 and 10 stay unchanged.

# Row 2
$ws.Range("D2").Value = 44320
$ws.Range("M2").Value = 80

# Row 3
$ws.Range("D3").Value = 44893
$ws.Range("M3").Value = 80
$ws.Range("N3").Value = 21000
$ws.Range("O3").Value = 22000
$ws.Range("P3").Value = 21625
$ws.Range("S3").Value = 1081

# Row 4
$ws.Range("D4").Value = 44533
$ws.Range("N4").Value = 16000
$ws.Range("O4").Value = 17000
$ws.Range("P4").Value = 16500
$ws.Range("S4").Value = 825

# Row 6
$ws.Range("D6").Value = 44708
$ws.Range("N6").Value = 20000
$ws.Range("O6").Value = 21000
$ws.Range("P6").Value = 20500
$ws.Range("S6").Value = 1025

# Row 7
$ws.Range("D7").Value = 44357
$ws.Range("M7").Value = 100
$ws.Range("N7").Value = 14000
$ws.Range("O7").Value = 15000
$ws.Range("P7").Value = 14500
$ws.Range("S7").Value = 725

# Row 8
$ws.Range("D8").Value = 44761
$ws.Range("M8").Value = 100
$ws.Range("N8").Value = 20000
$ws.Range("O8").Value = 21000
$ws.Range("P8").Value = 20500
$ws.Range("S8").Value = 1025
